{"js": "// Insert a new \"Author\" paragraph right after the \"Edison Achalma\" author\n// paragraph, containing the institutional affiliation text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Author\" && p.text.trim() === \"Edison Achalma\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Edison Achalma' Author paragraph\");\n}\n\n// NOTE: `Paragraph.insertParagraph(text, \"After\")` on this host moves the\n// *existing* paragraph's text into the freshly split paragraph (i.e. the\n// anchor keeps the empty half), so inserting text that way would corrupt\n// \"Edison Achalma\". Instead, insert a literal paragraph mark (\"\\r\") plus the\n// new text at the end of the target paragraph's range - this keeps \"Edison\n// Achalma\" in place and creates the new paragraph right after it, inheriting\n// the \"Author\" style from the split.\nconst endRange = target.getRange(\"End\");\nendRange.insertText(\n  \"\\rEscuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\",\n  \"After\"\n);\nawait context.sync();\n\n// Re-load paragraphs and make sure the newly created paragraph explicitly\n// carries the \"Author\" style (it should already inherit it, but set it\n// explicitly to be safe/robust).\nconst paragraphsAfter = body.paragraphs;\nparagraphsAfter.load(\"items,text,style\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphsAfter.items.length; i++) {\n  const p = paragraphsAfter.items[i];\n  if (\n    p.text.trim() ===\n    \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\"\n  ) {\n    p.style = \"Author\";\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Insert a new \"Author\" paragraph right after the \"Edison Achalma\" author\n# paragraph, containing the institutional affiliation text.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    $text = $p.Range.Text.Trim()\n    if ($styleName -eq \"Author\" -and $text -eq \"Edison Achalma\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Could not find the 'Edison Achalma' Author paragraph\"\n}\n\n# NOTE: collapsing to the end of the paragraph's range and then calling\n# InsertParagraphAfter() on it actually empties *this* paragraph and shifts\n# \"Edison Achalma\" into the freshly minted paragraph instead of leaving it in\n# place. Inserting the new text followed by a literal paragraph mark\n# (carriage return) via InsertAfter keeps \"Edison Achalma\" untouched and\n# creates the new paragraph right after it. (The new paragraph's style isn't\n# reliably inherited this way, so it is set explicitly below.)\n$r = $target.Range\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertAfter(\"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\" + [char]13)\n\n# Re-find the freshly created paragraph and make sure it explicitly carries\n# the \"Author\" style (it already inherits it from the split, but set it\n# explicitly to be safe/robust).\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.Trim()\n    if ($text -eq \"Escuela Profesional de Econom\u00eda, Universidad Nacional de San Crist\u00f3bal de Huamanga\") {\n        $p.Style = \"Author\"\n        break\n    }\n}\n"}
